# Concentraciones_Salida.xlsx -- refine the output-concentration table from 5
# to 10 time samples (halving dt) by inserting one new column before each
# existing "C(t = ...)" column and filling it with the newly simulated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank columns, left to right. Because each Insert() only shifts
# columns at-or-right-of the target, C/E/G/I/K are always still valid
# insertion points for the *next* original time column at the moment we use them.
foreach ($col in @("C", "E", "G", "I", "K")) {
    $ws.Columns($col).Insert()
}

# Header row for the 5 newly-inserted time columns
$ws.Range("C1").Value = "C(t = 281.250 s)"
$ws.Range("E1").Value = "C(t = 843.750 s)"
$ws.Range("G1").Value = "C(t = 1406.250 s)"
$ws.Range("I1").Value = "C(t = 1968.750 s)"
$ws.Range("K1").Value = "C(t = 2531.250 s)"

# Column C: newly simulated concentrations for rows 2-42
$Cvals = @(
    "0", "2.232367316034164e-09", "8.416902878718349e-09", "2.937991608498348e-08", "1.014001997939559e-07",
    "3.476057067870299e-07", "1.183245935430947e-06", "3.996327410193308e-06", "1.337876586823851e-05",
    "4.43418279662544e-05", "0.0001452773225432103", "0.000469606904874385", "0.001493946684086506",
    "0.004661597564137691", "0.01419982568061658", "0.04423999634441782", "0.121304769329773",
    "0.2720502201255362", "0.4823758125045884", "0.6790621757529187", "0.7598667626840945",
    "0.6790621757529187", "0.4823758125045884", "0.2720502201255361", "0.121304769329773",
    "0.04423999634441782", "0.01419982568061658", "0.004661597564137689", "0.001493946684086505",
    "0.0004696069048743845", "0.0001452773225432102", "4.434182796625434e-05", "1.337876586823849e-05",
    "3.996327410193301e-06", "1.183245935430945e-06", "3.476057067870293e-07", "1.014001997939557e-07",
    "2.937991608498342e-08", "8.416902878718331e-09", "2.232367316034159e-09", "0"
)
for ($i = 0; $i -lt $Cvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = [double]$Cvals[$i]
}

# Column E: newly simulated concentrations for rows 2-42
$Evals = @(
    "0", "3.450836868975576e-06", "1.04860752018988e-05", "2.80407620078825e-05", "7.241621754445021e-05",
    "0.0001821914338147642", "0.0004456855581129154", "0.001055832416525208", "0.002410068162060757",
    "0.00526980886144969", "0.01096958533241334", "0.02160621834085773", "0.04004942451908836",
    "0.06954706516333534", "0.1127456446648084", "0.1701918920722868", "0.2387892194202114",
    "0.3110339176759708", "0.3758281924485878", "0.4210803843345883", "0.437354794286582",
    "0.4210803843345883", "0.3758281924485878", "0.3110339176759707", "0.2387892194202114",
    "0.1701918920722867", "0.1127456446648083", "0.06954706516333531", "0.04004942451908834",
    "0.02160621834085771", "0.01096958533241334", "0.005269808861449685", "0.002410068162060754",
    "0.001055832416525207", "0.0004456855581129147", "0.000182191433814764", "7.24162175444501e-05",
    "2.804076200788245e-05", "1.048607520189878e-05", "3.450836868975569e-06", "0"
)
for ($i = 0; $i -lt $Evals.Length; $i++) {
    $row = $i + 2
    $ws.Range("E$row").Value = [double]$Evals[$i]
}

# Column G: newly simulated concentrations for rows 2-42
$Gvals = @(
    "0", "0.0001186390809159283", "0.000301193078222778", "0.0006365980149902988", "0.001268380461844555",
    "0.002431980453999468", "0.004499813258811713", "0.008028289970889568", "0.01379268625871305",
    "0.02278648389430042", "0.03615560548679386", "0.05504078485630945", "0.08031888977300591",
    "0.1122679290953847", "0.150224198926883", "0.1923376679590826", "0.2355421670562428",
    "0.2758235207359012", "0.3087904104941809", "0.3304504322644085", "0.3380058888268493",
    "0.3304504322644085", "0.3087904104941808", "0.2758235207359012", "0.2355421670562428",
    "0.1923376679590825", "0.1502241989268829", "0.1122679290953846", "0.08031888977300586",
    "0.05504078485630941", "0.03615560548679383", "0.0227864838943004", "0.01379268625871304",
    "0.008028289970889561", "0.004499813258811708", "0.002431980453999465", "0.001268380461844554",
    "0.0006365980149902978", "0.0003011930782227775", "0.0001186390809159281", "0"
)
for ($i = 0; $i -lt $Gvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = [double]$Gvals[$i]
}

# Column I: newly simulated concentrations for rows 2-42
$Ivals = @(
    "0", "0.000700913902878871", "0.001599901871726558", "0.002924585409562403", "0.00495946192902",
    "0.008067767955715743", "0.01270292747232854", "0.01940295517147655", "0.02876051335476145",
    "0.04136275072255449", "0.0576994262107754", "0.0780452699191895", "0.1023320362762567",
    "0.1300349418336317", "0.1601038861954337", "0.1909687527513076", "0.220638201752436",
    "0.2468932508158437", "0.2675541905840637", "0.2807780611095", "0.2853307716802928", "0.2807780611095",
    "0.2675541905840636", "0.2468932508158436", "0.2206382017524359", "0.1909687527513075",
    "0.1601038861954336", "0.1300349418336316", "0.1023320362762566", "0.07804526991918945",
    "0.05769942621077533", "0.04136275072255445", "0.02876051335476142", "0.01940295517147653",
    "0.01270292747232853", "0.008067767955715731", "0.004959461929019994", "0.002924585409562398",
    "0.001599901871726555", "0.0007009139028788698", "0"
)
for ($i = 0; $i -lt $Ivals.Length; $i++) {
    $row = $i + 2
    $ws.Range("I$row").Value = [double]$Ivals[$i]
}

# Column K: newly simulated concentrations for rows 2-42
$Kvals = @(
    "0", "0.001844089161787486", "0.003986733481069819", "0.006740457676998732", "0.01044186245959752",
    "0.0154541337198032", "0.02215861306322456", "0.03093296301068205", "0.04211502616882774",
    "0.05595381922056928", "0.07255207155341389", "0.09180787358880595", "0.1133656314305319",
    "0.1365877836199367", "0.1605578604796123", "0.1841220512817366", "0.2059706713562811",
    "0.2247536581809918", "0.2392169240534479", "0.2483407769128555", "0.2514592325607921",
    "0.2483407769128554", "0.2392169240534479", "0.2247536581809917", "0.2059706713562809",
    "0.1841220512817365", "0.1605578604796122", "0.1365877836199366", "0.1133656314305318",
    "0.09180787358880582", "0.07255207155341381", "0.05595381922056921", "0.04211502616882767",
    "0.03093296301068201", "0.02215861306322452", "0.01545413371980317", "0.0104418624595975",
    "0.00674045767699872", "0.003986733481069812", "0.001844089161787483", "0"
)
for ($i = 0; $i -lt $Kvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("K$row").Value = [double]$Kvals[$i]
}
